$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 920.4767212210036
$ws.Range("D11").Value = 920.4767212210036
$ws.Range("D12").Value = 830.5004175541554
$ws.Range("D13").Value = 830.5004175541554
